$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9016236867239733
$ws.Range("B3").Value = 0.9169381107491856
$ws.Range("B4").Value = 0.9154471544715447
$ws.Range("B5").Value = 0.9161920260374288
